# "lat to long and vice versa"
# Header row 1 gets three label corrections:
#   H1: Email          -> ProfilePicURL
#   I1: Longtitude      -> Latitude
#   J1: Latitude        -> Longitude
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H1").Value = "ProfilePicURL"
$ws.Range("I1").Value = "Latitude"
$ws.Range("J1").Value = "Longitude"
